$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("Part"), shifting Instruction..Picture from B..I to C..J
$ws.Columns("B:B").Insert()

# New header for column B
$ws.Range("B1").Value = "Part"

# Column A previously held "1a"/"1b"/"1c" combining question number + part letter.
# Split into A = question number (1) and B = part number (1,2,3)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 3

# Column widths: new column B should match column A (width 12, bestFit)
$ws.Columns("B:B").ColumnWidth = 11.17

# New column J (Picture, formerly I) width 10
$ws.Columns("J:J").ColumnWidth = 9.17

# Restore selection to A3 as in target sheetView
[void]$ws.Range("A3").Select()
